# Generate Report for Handoff
# Update the "Latest Handoff Date" / "Latest Handoff Datetime" for the
# "0a831913-0fb3-420a-8e24-486b73839eae" entry after a new handoff was
# generated.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D5").Value = "2016-03-23 20:44:28"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E5").Value = "2016-03-23 20:44:24"
